# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp banner in A1
$ws.Range("A1").Value = "Datos actualizados a 21 de Mayo de 2020 a las 06:35"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1593039
$ws.Range("C4").Value = 316
$ws.Range("D4").Value = 370812
$ws.Range("E4").Value = 1127286
$ws.Range("G4").Value = 5
$ws.Range("H4").Value = 94941

# Row 14 - India
$ws.Range("B14").Value = 112442
$ws.Range("C14").Value = 414
$ws.Range("E14").Value = 63582
$ws.Range("G14").Value = 4
$ws.Range("H14").Value = 3438

# Row 100 - Kirguistan
$ws.Range("B100").Value = 1313
$ws.Range("C100").Value = 43
$ws.Range("D100").Value = 923
$ws.Range("E100").Value = 376

# Row 188 - Botsuana
$ws.Range("B188").Value = 29
$ws.Range("C188").Value = 4
$ws.Range("E188").Value = 11

# Row 193 - Butan
$ws.Range("D193").Value = 6
$ws.Range("E193").Value = 15
